$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 21483.229
$ws.Range("I15").Value = 21483.229
$ws.Range("K15").Value = 64449.687
$ws.Range("M15").Value = -64280.687
$ws.Range("H80").Value = 11447162
$ws.Range("I80").Value = 20853000
$ws.Range("J80").Value = 160155.8
$ws.Range("K80").Value = 62559000
$ws.Range("L80").Value = 480467.4
$ws.Range("M80").Value = -62558002
$ws.Range("N80").Value = -482463.4
$ws.Range("H83").Value = 11447162
$ws.Range("I83").Value = 20853000
$ws.Range("J83").Value = 160155.8
$ws.Range("K83").Value = 187677000
$ws.Range("L83").Value = 1441402.2
$ws.Range("M83").Value = -187672008
$ws.Range("N83").Value = -1451386.2
$ws.Range("H86").Value = 198416400
$ws.Range("I86").Value = 266667860
$ws.Range("K86").Value = 266667860
$ws.Range("M86").Value = -266666737
$ws.Range("H89").Value = 198416400
$ws.Range("I89").Value = 266667860
$ws.Range("K89").Value = 1333339300
$ws.Range("M89").Value = -1333333684
$ws.Range("H109").Value = 43849.5
$ws.Range("J109").Value = 43849.5
$ws.Range("L109").Value = 43849.5
$ws.Range("N109").Value = -46623.5
$ws.Range("H138").Value = 3337474
$ws.Range("I138").Value = 2012.7222
$ws.Range("J138").Value = 8340666
$ws.Range("K138").Value = 6038.1666
$ws.Range("L138").Value = 25021998
$ws.Range("M138").Value = -898.1665999999996
$ws.Range("N138").Value = -25032278

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2362075
$ws.Range("I32").Value = 2553755
$ws.Range("K32").Value = 2553755
$ws.Range("M32").Value = -2553468
$ws.Range("H102").Value = 988.6667
$ws.Range("I102").Value = 924.75
$ws.Range("K102").Value = 924.75
$ws.Range("M102").Value = 697.25
$ws.Range("H110").Value = 55557200
$ws.Range("I110").Value = 1462.5
$ws.Range("K110").Value = 1462.5
$ws.Range("M110").Value = 582.5
$ws.Range("H132").Value = 1860556.4
$ws.Range("I132").Value = 5009969
$ws.Range("J132").Value = 7960.7646
$ws.Range("K132").Value = 15029907
$ws.Range("L132").Value = 23882.2938
$ws.Range("M132").Value = -15027377
$ws.Range("N132").Value = -28942.2938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 16668726
$ws.Range("J20").Value = 6816.5
$ws.Range("K20").Value = 16668726
$ws.Range("L20").Value = 6816.5
$ws.Range("M20").Value = -16668479
$ws.Range("N20").Value = -7310.5
$ws.Range("H86").Value = 9023998
$ws.Range("I86").Value = 17944868
$ws.Range("K86").Value = 17944868
$ws.Range("M86").Value = -17943745
$ws.Range("H89").Value = 9023998
$ws.Range("I89").Value = 17944868
$ws.Range("K89").Value = 89724340
$ws.Range("M89").Value = -89718724
$ws.Range("H94").Value = 41670172
$ws.Range("J94").Value = 7008
$ws.Range("L94").Value = 7008
$ws.Range("N94").Value = -7910
$ws.Range("H134").Value = 7904.2607
$ws.Range("I134").Value = 3045.75
$ws.Range("K134").Value = 9137.25
$ws.Range("M134").Value = -6602.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 2000
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 200.33333
$ws.Range("I31").Value = 250.5
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 751.5
$ws.Range("L31").Value = 300
$ws.Range("M31").Value = -463.5
$ws.Range("N31").Value = -876
$ws.Range("H39").Value = 9215.799999999999
$ws.Range("I39").Value = 400
$ws.Range("J39").Value = 9583.125
$ws.Range("K39").Value = 1200
$ws.Range("L39").Value = 28749.375
$ws.Range("M39").Value = -906
$ws.Range("N39").Value = -29337.375
$ws.Range("H68").Value = 2193.7585
$ws.Range("J68").Value = 2427.4285
$ws.Range("L68").Value = 7282.2855
$ws.Range("N68").Value = -8904.2855
$ws.Range("H71").Value = 2193.7585
$ws.Range("J71").Value = 2427.4285
$ws.Range("L71").Value = 21846.8565
$ws.Range("N71").Value = -29958.8565
$ws.Range("H81").Value = 7137.6665
$ws.Range("J81").Value = 7137.6665
$ws.Range("L81").Value = 21412.9995
$ws.Range("N81").Value = -23658.9995
$ws.Range("H84").Value = 7137.6665
$ws.Range("J84").Value = 7137.6665
$ws.Range("L84").Value = 64238.9985
$ws.Range("N84").Value = -75470.9985
$ws.Range("H86").Value = 3666.6667
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 9000
$ws.Range("L86").Value = 15000
$ws.Range("M86").Value = -7814
$ws.Range("N86").Value = -17372
$ws.Range("H89").Value = 3666.6667
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 27000
$ws.Range("L89").Value = 45000
$ws.Range("M89").Value = -21072
$ws.Range("N89").Value = -56856
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").Value = 0
$ws.Range("H125").Value = 5750
$ws.Range("I125").Value = 5700
$ws.Range("J125").Value = 5800
$ws.Range("K125").Value = 17100
$ws.Range("L125").Value = 17400
$ws.Range("M125").Value = -12180
$ws.Range("N125").Value = -27240
$ws.Range("H132").Value = 8116.227
$ws.Range("I132").Value = 2180.0908
$ws.Range("J132").Value = 14052.363
$ws.Range("K132").Value = 19620.8172
$ws.Range("L132").Value = 126471.267
$ws.Range("M132").Value = -17090.8172
$ws.Range("N132").Value = -131531.267

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 202899.8
$ws.Range("J80").Value = 202899.8
$ws.Range("L80").Value = 202899.8
$ws.Range("N80").Value = -204895.8
$ws.Range("H83").Value = 202899.8
$ws.Range("J83").Value = 202899.8
$ws.Range("L83").Value = 1014499
$ws.Range("N83").Value = -1024483
$ws.Range("H102").Value = 4059.3333
$ws.Range("I102").Value = 4158.28
$ws.Range("K102").Value = 4158.28
$ws.Range("M102").Value = -2536.28
$ws.Range("H113").Value = 6703.4243
$ws.Range("I113").Value = 4141.6665
$ws.Range("J113").Value = 8167.2856
$ws.Range("K113").Value = 4141.6665
$ws.Range("L113").Value = 8167.2856
$ws.Range("M113").Value = -1971.6665
$ws.Range("N113").Value = -12507.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 984.875
$ws.Range("I55").Value = 1091.2142
$ws.Range("J55").Value = 836
$ws.Range("K55").Value = 1091.2142
$ws.Range("L55").Value = 836
$ws.Range("M55").Value = -918.2141999999999
$ws.Range("N55").Value = -1182
$ws.Range("H122").Value = 8016.9165
$ws.Range("I122").Value = 8347.764999999999
$ws.Range("J122").Value = 7213.4287
$ws.Range("K122").Value = 25043.295
$ws.Range("L122").Value = 21640.2861
$ws.Range("M122").Value = -22593.295
$ws.Range("N122").Value = -26540.2861
$ws.Range("H132").Value = 6706.486
$ws.Range("I132").Value = 4106.4
$ws.Range("K132").Value = 12319.2
$ws.Range("M132").Value = -9789.199999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3041.8333
$ws.Range("I100").Value = 2834.2222
$ws.Range("K100").Value = 5668.4444
$ws.Range("M100").Value = -5127.4444
$ws.Range("H107").Value = 910.3333
$ws.Range("J107").Value = 720
$ws.Range("L107").Value = 2160
$ws.Range("N107").Value = -6000
